# Update stand and report
# Updates the "informe interaccion stand" worksheet: refreshes the
# interaction-with-stand-elements table (rows 3-13) with new data,
# and removes the now-obsolete last row (old row 14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("informe interaccion stand")

# --- Row 3 ---------------------------------------------------------
$ws.Cells.Item(3, 3).Value  = "Stand 5"            # C3
$ws.Cells.Item(3, 4).Value  = 1                    # D3
$ws.Cells.Item(3, 6).Value  = 44230.4423611111     # F3
$ws.Cells.Item(3, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(3, 7).Value  = "Goberto Calleja"    # G3
$ws.Cells.Item(3, 8).Value  = "Secpho"             # H3
$ws.Cells.Item(3, 9).Value  = "Engineering"        # I3
$ws.Cells.Item(3, 10).Value = "CEO"                # J3
$ws.Cells.Item(3, 11).Value = "Stand 5"            # K3
$ws.Cells.Item(3, 13).Value = "Goberto Calleja"    # M3
$ws.Cells.Item(3, 14).Value = "Secpho"             # N3
$ws.Cells.Item(3, 15).Value = "Engineering"        # O3
$ws.Cells.Item(3, 16).Value = "CEO"                # P3
$ws.Cells.Item(3, 17).Value = "user1@gmail.com"    # Q3

# --- Row 4 ---------------------------------------------------------
$ws.Cells.Item(4, 3).Value  = "Stand 6"            # C4
$ws.Cells.Item(4, 4).Value  = 7                    # D4
$ws.Cells.Item(4, 6).Value  = 44230.4625           # F4
$ws.Cells.Item(4, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(4, 7).Value  = "Goberto Calleja"    # G4
$ws.Cells.Item(4, 8).Value  = "Secpho"             # H4
$ws.Cells.Item(4, 9).Value  = "Engineering"        # I4
$ws.Cells.Item(4, 10).Value = "CEO"                # J4
$ws.Cells.Item(4, 11).Value = "Stand 6"            # K4
# old M4:O4 + Q4 values no longer apply to this row
$ws.Cells.Item(4, 13).ClearContents()
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(4, 15).ClearContents()
$ws.Cells.Item(4, 17).ClearContents()

# --- Row 5 ---------------------------------------------------------
$ws.Cells.Item(5, 3).Value  = "AsorCAD counter"    # C5
$ws.Cells.Item(5, 4).Value  = 3                    # D5
$ws.Cells.Item(5, 6).Value  = 44230.4673611111     # F5
$ws.Cells.Item(5, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(5, 7).Value  = "Goberto Calleja"    # G5
$ws.Cells.Item(5, 8).Value  = "Secpho"             # H5
$ws.Cells.Item(5, 9).Value  = "Engineering"        # I5
$ws.Cells.Item(5, 10).Value = "CEO"                # J5
$ws.Cells.Item(5, 11).Value = "Stand 6"            # K5
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(5, 15).ClearContents()
$ws.Cells.Item(5, 17).ClearContents()

# --- Row 6 (the stand-counter columns C/D no longer used here) -----
$ws.Cells.Item(6, 3).ClearContents()
$ws.Cells.Item(6, 4).ClearContents()
$ws.Cells.Item(6, 6).Value  = 44230.6              # F6
$ws.Cells.Item(6, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(6, 7).Value  = "Pedro Carrillo"     # G6
$ws.Cells.Item(6, 8).Value  = "Secpho"             # H6
$ws.Cells.Item(6, 9).Value  = "Engineering"        # I6
$ws.Cells.Item(6, 10).Value = "prueba"             # J6
$ws.Cells.Item(6, 11).Value = "Stand 6"            # K6

# --- Row 7 -----------------------------------------------------------
$ws.Cells.Item(7, 6).Value  = 44230.6              # F7
$ws.Cells.Item(7, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(7, 7).Value  = "Pedro Carrillo"     # G7
$ws.Cells.Item(7, 8).Value  = "Secpho"             # H7
$ws.Cells.Item(7, 9).Value  = "Engineering"        # I7
$ws.Cells.Item(7, 10).Value = "prueba"             # J7
$ws.Cells.Item(7, 11).Value = "Stand 6"            # K7

# --- Row 8 -----------------------------------------------------------
$ws.Cells.Item(8, 6).Value  = 44230.6              # F8
$ws.Cells.Item(8, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(8, 7).Value  = "Pedro Carrillo"     # G8
$ws.Cells.Item(8, 8).Value  = "Secpho"             # H8
$ws.Cells.Item(8, 9).Value  = "Engineering"        # I8
$ws.Cells.Item(8, 10).Value = "prueba"             # J8
$ws.Cells.Item(8, 11).Value = "Stand 6"            # K8

# --- Row 9 -----------------------------------------------------------
$ws.Cells.Item(9, 6).Value  = 44230.6              # F9
$ws.Cells.Item(9, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(9, 7).Value  = "Juan Camilla"       # G9
$ws.Cells.Item(9, 8).Value  = "Secpho"             # H9
$ws.Cells.Item(9, 9).Value  = "Engineering"        # I9
$ws.Cells.Item(9, 10).Value = "prueba"             # J9
$ws.Cells.Item(9, 11).Value = "Stand 6"            # K9

# --- Row 10 ----------------------------------------------------------
$ws.Cells.Item(10, 6).Value  = 44230.6             # F10
$ws.Cells.Item(10, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(10, 7).Value  = "Juan Camilla"      # G10
$ws.Cells.Item(10, 8).Value  = "Secpho"            # H10
$ws.Cells.Item(10, 9).Value  = "Engineering"       # I10
$ws.Cells.Item(10, 10).Value = "prueba"            # J10
$ws.Cells.Item(10, 11).Value = "Stand 6"           # K10

# --- Row 11 ------------------------------------------------------------
$ws.Cells.Item(11, 6).Value  = 44232.5423611111    # F11
$ws.Cells.Item(11, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(11, 7).Value  = "Goberto Calleja"   # G11
$ws.Cells.Item(11, 8).Value  = "Secpho"            # H11
$ws.Cells.Item(11, 9).Value  = "Engineering"       # I11
$ws.Cells.Item(11, 10).Value = "CEO"               # J11
$ws.Cells.Item(11, 11).Value = "AsorCAD counter"   # K11

# --- Row 12 --------------------------------------------------------------
$ws.Cells.Item(12, 6).Value  = 44232.5423611111    # F12
$ws.Cells.Item(12, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(12, 7).Value  = "Goberto Calleja"   # G12
$ws.Cells.Item(12, 8).Value  = "Secpho"            # H12
$ws.Cells.Item(12, 9).Value  = "Engineering"       # I12
$ws.Cells.Item(12, 10).Value = "CEO"               # J12
$ws.Cells.Item(12, 11).Value = "AsorCAD counter"   # K12

# --- Row 13 --------------------------------------------------------------
$ws.Cells.Item(13, 6).Value  = 44232.5444444444    # F13
$ws.Cells.Item(13, 6).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(13, 7).Value  = "Goberto Calleja"   # G13
$ws.Cells.Item(13, 8).Value  = "Secpho"            # H13
$ws.Cells.Item(13, 9).Value  = "Engineering"       # I13
$ws.Cells.Item(13, 10).Value = "CEO"               # J13
$ws.Cells.Item(13, 11).Value = "AsorCAD counter"   # K13

# Old row 14 is gone entirely in the updated report -> delete it,
# shrinking the sheet dimension from A1:Q14 down to A1:Q13.
$ws.Rows("14:14").Delete()
